# Apply the Tbl_Admin_PedMedCont sheet updates described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tbl_Admin_PedMedCont")

# --- Row 14 ---
$ws.Range("P14").Value = 5

# --- Row 17 ---
$ws.Range("D17").Value = 1.5
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 2.5
$ws.Range("G17").Value = 50
$ws.Range("H17").Value = 5
$ws.Range("I17").Value = 50
$ws.Range("J17").Value = 5
$ws.Range("K17").Value = 50
$ws.Range("P17").Value = 10

# --- Row 19 ---
$ws.Range("P19").Value = 1000

# --- Row 23 ---
$ws.Range("O23").Value = 1
$ws.Range("P23").Value = 16

# --- Row 27 ---
$ws.Range("Q27").Value = 1

# --- Row 33 ---
$ws.Range("O33").Value = 0.005
$ws.Range("P33").Value = 0.42

# --- Row 34 ---
$ws.Range("D34").Value = 1250
$ws.Range("F34").Value = 1250
$ws.Range("H34").Value = 1250
$ws.Range("J34").Value = 1250
$ws.Range("O34").Value = 5
$ws.Range("P34").Value = 25

# --- Selection moves to B8 ---
$ws.Range("B8").Select()
